$wb = $excel.ActiveWorkbook

# --- "Prix Spot": a new day ("07-nov") of (still unpublished) hourly prices
#     is inserted as a column just before the October block, shifting every
#     column from the old "01-oct." onward one place to the right. ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("DL1").EntireColumn.Insert()
$ws1.Range("DL1").Value = "07-nov"
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 116).Value = "-"
}

# --- "Gaz": append the next day's price row ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A144").NumberFormat = "@"
$ws2.Range("A144").Value = "2025-11-05"
$ws2.Range("A144").Style = "Normal"
$ws2.Range("B144").Value = 30.425

# --- "CO2": append the next day's price row ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A144").NumberFormat = "@"
$ws3.Range("A144").Value = "2025-11-05"
$ws3.Range("A144").Style = "Normal"
$ws3.Range("B144").Value = 81.18
